# Updated symbol list (coin prices / 1h volume %) on Sun Jan 29 22:36:26 UTC 2023
# with GitHub Actions. Also a few coin rows (8-17) shifted by one position
# (new GateToken row inserted, rest pushed down).
#
# All Price/Volume cells in this sheet are stored as literal text (not
# numbers), e.g. D2 = "317.92" and E2 = "4.01%". Plain `Range.Value = "..."`
# assignment lets Excel auto-coerce numeric-looking strings into real
# numbers (and silently applies a Number/Percent style), which would change
# both the cell type and its style away from the original "General"/no-style
# cells. Set-TextValue forces the text interpretation via a temporary "@"
# (Text) number format, then strips the format back off with ClearFormats()
# so the cell ends up with no explicit style again - same as the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

Set-TextValue 'D2' '317.92'
Set-TextValue 'E2' '4.01%'
Set-TextValue 'D3' '39.82'
Set-TextValue 'E3' '2.55%'
Set-TextValue 'D4' '5.135'
Set-TextValue 'E4' '0.44%'
Set-TextValue 'D5' '0.08217'
Set-TextValue 'E5' '1.87%'
Set-TextValue 'D6' '2.062'
Set-TextValue 'E6' '6.68%'
Set-TextValue 'D7' '8.381'
Set-TextValue 'E7' '4.41%'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D8' '4.317'
Set-TextValue 'E8' '2.59%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9426'
Set-TextValue 'E9' '1.71%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1356'
Set-TextValue 'E10' '-6.08%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.2002'
Set-TextValue 'E11' '4.42%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.09093'
Set-TextValue 'E12' '0.60%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03526'
Set-TextValue 'E13' '0.57%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09798'
Set-TextValue 'E14' '0.24%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001408'
Set-TextValue 'E15' '0.49%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.006163'
Set-TextValue 'E16' '4.52%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.683'
Set-TextValue 'E17' '-2.00%'
Set-TextValue 'E18' '-0.59%'
Set-TextValue 'D19' '0.3494'
Set-TextValue 'E19' '0.93%'
Set-TextValue 'D20' '0.1322'
Set-TextValue 'E20' '-0.41%'
Set-TextValue 'D21' '4.949'
Set-TextValue 'E21' '5.40%'
Set-TextValue 'D23' '0.04367'
Set-TextValue 'E23' '-0.23%'
Set-TextValue 'D24' '0.001229'
Set-TextValue 'E24' '1.76%'
Set-TextValue 'D25' '0.004794'
Set-TextValue 'E25' '12.35%'
Set-TextValue 'E26' '-0.03%'
Set-TextValue 'D27' '0.0004000'
Set-TextValue 'E27' '-10.06%'
Set-TextValue 'D39' '0.02326'
Set-TextValue 'E39' '14.40%'
Set-TextValue 'D40' '0.05196'
Set-TextValue 'E40' '2.98%'
Set-TextValue 'D41' '0.007749'
Set-TextValue 'E41' '3.09%'
Set-TextValue 'D42' '0.009872'
Set-TextValue 'E42' '0.97%'
Set-TextValue 'E43' '4.83%'
Set-TextValue 'D44' '0.002081'
Set-TextValue 'E44' '-0.61%'
Set-TextValue 'D45' '0.008910'
Set-TextValue 'E45' '-9.38%'
Set-TextValue 'D46' '0.00006611'
Set-TextValue 'E46' '6.44%'
Set-TextValue 'E47' '-0.16%'
Set-TextValue 'D48' '0.002945'
Set-TextValue 'E48' '2.52%'
Set-TextValue 'E49' '-6.29%'
Set-TextValue 'E50' '-0.16%'
Set-TextValue 'E51' '-0.16%'
